# The captured diff for this revision is a pure OOXML re-serialization:
# every hunk in word/document.xml (sectPr), word/numbering.xml and
# word/styles.xml only reorders existing XML attributes (e.g.
# w:abstractNumId/w15:restartNumberingAfterBreak, w:val/w:pos on <w:tab>,
# w:left/w:hanging on <w:ind>, w:val/w:themeColor/w:themeShade on
# <w:color>, etc.) - the attribute *values* are unchanged, and none of
# this content is reachable through the Word object model (attribute
# order is not something Find/Replace, Styles, PageSetup, or
# ListTemplates expose or let you control).
#
# The single literal value change in the diff - the <w:nsid> GUID on the
# orphaned abstractNum/num definition (abstractNumId="990" / numId
# "1000") - is also not reachable via COM: Word does not expose NSID on
# List/ListTemplate/ListLevel, and that list definition is not applied
# to any paragraph in the document (no w:numPr anywhere in
# word/document.xml), so it has no visible/semantic effect on the
# document.
#
# Net effect: there is no content, formatting, or layout change for
# Word automation to perform here. Touch the document via the object
# model without altering anything, so the package is re-saved cleanly
# and no spurious edits are introduced.
$d = $word.ActiveDocument
$null = $d.Content
